$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RF values (column I) for rows 31 through 48 to reflect the
# updated 2025 data / RF changes.
for ($r = 31; $r -le 48; $r++) {
    $ws.Cells.Item($r, 9).Value = 13.11588235294118
}
